$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = ''''
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 37
$ws.Range("D8").Value = '''1.0'
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = '''0.00'
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = ''''

# Row 9
$ws.Range("A9").Value = ''''
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 78
$ws.Range("D9").Value = '''2.0'
$ws.Range("E9").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = '''0.00'
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = ''''

# Row 10
$ws.Range("A10").Value = 'P. point'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 91
$ws.Range("D10").Value = '''6'
$ws.Range("E10").Value = 'On board'
$ws.Range("F10").Value = 136
$ws.Range("G10").Value = '''12376.00'
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = ''''

# Row 11
$ws.Range("A11").Value = 'Each'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 22
$ws.Range("D11").Value = '''4.0'
$ws.Range("E11").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = '''1100.00'
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = ''''

# Row 12
$ws.Range("A12").Value = ''''
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 82
$ws.Range("D12").Value = '''11.0'
$ws.Range("E12").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = '''0.00'
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = ''''

# Row 13
$ws.Range("A13").Value = 'R. mtr.'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 87
$ws.Range("D13").Value = '''17'
$ws.Range("E13").Value = '25 mm'
$ws.Range("F13").Value = 56
$ws.Range("G13").Value = '''4872.00'
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = ''''

# Row 14
$ws.Range("A14").Value = 'Mtr.'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 92
$ws.Range("D14").Value = '''20'
$ws.Range("E14").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F14").Value = 122
$ws.Range("G14").Value = '''11224.00'
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = ''''

# Row 15
$ws.Range("A15").Value = 'Mtr.'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 22
$ws.Range("D15").Value = '''23'
$ws.Range("E15").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = '''440.00'
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = ''''

# Row 16
$ws.Range("A16").Value = ''''
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = '''16.0'
$ws.Range("E16").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = '''0.00'
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = ''''

# Row 17
$ws.Range("A17").Value = 'Each'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = '''30'
$ws.Range("E17").Value = ' 6 A to 32 A rating'
$ws.Range("F17").Value = 187
$ws.Range("G17").Value = '''561.00'
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = ''''

# Row 18
$ws.Range("A18").Value = ''''
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 63
$ws.Range("D18").Value = '''18.0'
$ws.Range("E18").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = '''0.00'
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = ''''

# Row 19
$ws.Range("A19").Value = ''''
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 70
$ws.Range("D19").Value = '''36'
$ws.Range("E19").Value = 'Total'
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = '''0.00'
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = ''''

# Row 20
$ws.Range("A20").Value = ''''
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 81
$ws.Range("D20").Value = '''38'
$ws.Range("E20").Value = 'Grand Total'
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = '''0.00'
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = ''''

# Row 21
$ws.Range("A21").Value = ''''

# Row 22
$ws.Range("A22").Value = ''''
$ws.Range("B22").Value = ''''
$ws.Range("C22").Value = ''''
$ws.Range("D22").Value = ''''
$ws.Range("E22").Value = 'Grand Total Rs.'
$ws.Range("F22").Value = ''''
$ws.Range("G22").Value = '''30573.00'
$ws.Range("H22").Value = '''30573.00'
$ws.Range("I22").Value = ''''

# Row 23
$ws.Range("A23").Value = ''''
$ws.Range("B23").Value = ''''
$ws.Range("C23").Value = ''''
$ws.Range("D23").Value = ''''
$ws.Range("E23").Value = 'Tender Premium @ 0%'
$ws.Range("F23").Value = ''''
$ws.Range("G23").Value = '''0.00'
$ws.Range("H23").Value = '''0.00'
$ws.Range("I23").Value = ''''

# Row 24
$ws.Range("A24").Value = ''''
$ws.Range("B24").Value = ''''
$ws.Range("C24").Value = ''''
$ws.Range("D24").Value = ''''
$ws.Range("E24").Value = 'NET PAYABLE AMOUNT Rs.'
$ws.Range("F24").Value = ''''
$ws.Range("G24").Value = '''30573.00'
$ws.Range("H24").Value = '''30573.00'
$ws.Range("I24").Value = ''''

# Clear B21:I21 - only A21 remains with empty string
$ws.Range("B21:I21").ClearContents()
